$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Spent" (AC) actuals for the month grew from 200,000 to 300,000 -
# this is the substantive budget-tracker update; CV, CPI, EAC and VAC
# (F8:F11) are formulas and recalculate automatically.
$ws.Range("F6").Value = 300000

# The author's active cell moved to F7 by the time the file was saved.
$ws.Range("F7").Select()
